# Natmi following Dr Hou advice
#
# Recomputes the Penk -> Ogfr ligand-receptor pair stats (NATMI lrc2p output)
# after correcting the per-cluster cell counts, and adds the new "sCs"
# sending-cluster rows (rows 8-10) that come from the corrected cluster set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Penk"
$ws.Range("C2").Value = "Ogfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3329863333333333
$ws.Range("H2").Value = 0.9989589999999999
$ws.Range("I2").Value = 0.002397412224472272
$ws.Range("J2").Value = 0.002397412224472272
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.91824966666667
$ws.Range("N2").Value = 38.754749
$ws.Range("O2").Value = 0.5690183414937939
$ws.Range("P2").Value = 0.5690183414937939
$ws.Range("Q2").Value = 4.301600589587889
$ws.Range("R2").Value = 38.714405306291
$ws.Range("S2").Value = 0.00136417152784616
$ws.Range("T2").Value = 0.001364171527846159

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Penk"
$ws.Range("C3").Value = "Ogfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3329863333333333
$ws.Range("H3").Value = 0.9989589999999999
$ws.Range("I3").Value = 0.002397412224472272
$ws.Range("J3").Value = 0.002397412224472272
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.587813999999999
$ws.Range("N3").Value = 19.763442
$ws.Range("O3").Value = 0.2901776241422384
$ws.Range("P3").Value = 0.2901776241422384
$ws.Range("Q3").Value = 2.193652028542
$ws.Range("R3").Value = 19.742868256878
$ws.Range("S3").Value = 0.0006956753833869226
$ws.Range("T3").Value = 0.0006956753833869225

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Penk"
$ws.Range("C4").Value = "Ogfr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3329863333333333
$ws.Range("H4").Value = 0.9989589999999999
$ws.Range("I4").Value = 0.002397412224472272
$ws.Range("J4").Value = 0.002397412224472272
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.196631
$ws.Range("N4").Value = 9.589893
$ws.Range("O4").Value = 0.1408040343639677
$ws.Range("P4").Value = 0.1408040343639677
$ws.Range("Q4").Value = 1.064434435709667
$ws.Range("R4").Value = 9.579909921387
$ws.Range("S4").Value = 0.0003375653132391901
$ws.Range("T4").Value = 0.00033756531323919

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Penk"
$ws.Range("C5").Value = "Ogfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 138.2516836666667
$ws.Range("H5").Value = 414.755051
$ws.Range("I5").Value = 0.9953750148194476
$ws.Range("J5").Value = 0.9953750148194476
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.91824966666667
$ws.Range("N5").Value = 38.754749
$ws.Range("O5").Value = 0.5690183414937939
$ws.Range("P5").Value = 0.5690183414937939
$ws.Range("Q5").Value = 1785.969766443022
$ws.Range("R5").Value = 16073.7278979872
$ws.Range("S5").Value = 0.5663866400969226
$ws.Range("T5").Value = 0.5663866400969226

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Penk"
$ws.Range("C6").Value = "Ogfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 138.2516836666667
$ws.Range("H6").Value = 414.755051
$ws.Range("I6").Value = 0.9953750148194476
$ws.Range("J6").Value = 0.9953750148194476
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.587813999999999
$ws.Range("N6").Value = 19.763442
$ws.Range("O6").Value = 0.2901776241422384
$ws.Range("P6").Value = 0.2901776241422384
$ws.Range("Q6").Value = 910.7763771828378
$ws.Range("R6").Value = 8196.987394645541
$ws.Range("S6").Value = 0.2888355569308526
$ws.Range("T6").Value = 0.2888355569308526

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Penk"
$ws.Range("C7").Value = "Ogfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 138.2516836666667
$ws.Range("H7").Value = 414.755051
$ws.Range("I7").Value = 0.9953750148194476
$ws.Range("J7").Value = 0.9953750148194476
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.196631
$ws.Range("N7").Value = 9.589893
$ws.Range("O7").Value = 0.1408040343639677
$ws.Range("P7").Value = 0.1408040343639677
$ws.Range("Q7").Value = 441.9396178110603
$ws.Range("R7").Value = 3977.456560299543
$ws.Range("S7").Value = 0.1401528177916724
$ws.Range("T7").Value = 0.1401528177916724

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Penk"
$ws.Range("C8").Value = "Ogfr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3093966666666667
$ws.Range("H8").Value = 0.92819
$ws.Range("I8").Value = 0.002227572956080197
$ws.Range("J8").Value = 0.002227572956080197
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.91824966666667
$ws.Range("N8").Value = 38.754749
$ws.Range("O8").Value = 0.5690183414937939
$ws.Range("P8").Value = 0.5690183414937939
$ws.Range("Q8").Value = 3.996863386034445
$ws.Range("R8").Value = 35.97177047431
$ws.Range("S8").Value = 0.001267529869025182
$ws.Range("T8").Value = 0.001267529869025182

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Penk"
$ws.Range("C9").Value = "Ogfr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3093966666666667
$ws.Range("H9").Value = 0.92819
$ws.Range("I9").Value = 0.002227572956080197
$ws.Range("J9").Value = 0.002227572956080197
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.587813999999999
$ws.Range("N9").Value = 19.763442
$ws.Range("O9").Value = 0.2901776241422384
$ws.Range("P9").Value = 0.2901776241422384
$ws.Range("Q9").Value = 2.03824769222
$ws.Range("R9").Value = 18.34422922998
$ws.Range("S9").Value = 0.0006463918279988543
$ws.Range("T9").Value = 0.0006463918279988543

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Penk"
$ws.Range("C10").Value = "Ogfr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.3093966666666667
$ws.Range("H10").Value = 0.92819
$ws.Range("I10").Value = 0.002227572956080197
$ws.Range("J10").Value = 0.002227572956080197
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.196631
$ws.Range("N10").Value = 9.589893
$ws.Range("O10").Value = 0.1408040343639677
$ws.Range("P10").Value = 0.1408040343639677
$ws.Range("Q10").Value = 0.9890269759633333
$ws.Range("R10").Value = 8.90124278367
$ws.Range("S10").Value = 0.0003136512590561613
$ws.Range("T10").Value = 0.0003136512590561613

